$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "2"
$ws.Range("E2").Value = "47.50"

$ws.Range("C3").Value = "1"
$ws.Range("E3").Value = "13.00"

$ws.Range("C4").Value = "2"
$ws.Range("E4").Value = "26.00"

$ws.Range("C8").Value = "2"
$ws.Range("E8").Value = "18.50"

$ws.Range("C9").Value = "2"
$ws.Range("E9").Value = "18.50"

$ws.Range("C11").Value = "1"
$ws.Range("E11").Value = "14.00"
